# "Fruta / hortaliza, semanal" — add this week's price records for
# Pepino ensalada @ Terminal La Palmera de La Serena.
#
# A new week's pair of rows (Primera / Segunda quality) is inserted right
# before the existing "Región de Arica y Parinacota" block, pushing all
# subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 312:313 — everything from the old row 312
# onward shifts down to 314 onward.
$ws.Rows("312:313").Insert()

# New row 312 — Primera
$ws.Range("A312").Value2 = 8
$ws.Range("B312").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C312").Value2 = "Coquimbo"
$ws.Range("D312").Value2 = 44448
$ws.Range("E312").Value2 = 4
$ws.Range("F312").Value2 = 100112043
$ws.Range("G312").Value2 = "Pepino ensalada"
$ws.Range("H312").Value2 = "Sin especificar"
$ws.Range("I312").Value2 = "Primera"
$ws.Range("J312").Value2 = 600
$ws.Range("K312").Value2 = 15000
$ws.Range("L312").Value2 = 16000
$ws.Range("M312").Value2 = 15500
$ws.Range("N312").Value2 = "$/caja 60 unidades"
$ws.Range("O312").Value2 = "Región de Arica y Parinacota"
$ws.Range("P312").Value2 = 258
$ws.Range("Q312").Value2 = 60
$ws.Range("R312").Value2 = "Hortaliza"

# New row 313 — Segunda
$ws.Range("A313").Value2 = 8
$ws.Range("B313").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C313").Value2 = "Coquimbo"
$ws.Range("D313").Value2 = 44448
$ws.Range("E313").Value2 = 4
$ws.Range("F313").Value2 = 100112043
$ws.Range("G313").Value2 = "Pepino ensalada"
$ws.Range("H313").Value2 = "Sin especificar"
$ws.Range("I313").Value2 = "Segunda"
$ws.Range("J313").Value2 = 400
$ws.Range("K313").Value2 = 10000
$ws.Range("L313").Value2 = 11000
$ws.Range("M313").Value2 = 10500
$ws.Range("N313").Value2 = "$/caja 100 unidades"
$ws.Range("O313").Value2 = "Región de Arica y Parinacota"
$ws.Range("P313").Value2 = 105
$ws.Range("Q313").Value2 = 100
$ws.Range("R313").Value2 = "Hortaliza"
